# "New Update: Autocomplete input"
# - Update the item description / part-number / product name on the
#   "nhap-linhkien" sheet to the new DCU part (row 2), and correct the
#   quantity in stock.
# - Remove the now-obsolete export record (row 2) on the "xuat-linhkien"
#   sheet entirely.

$wb = $excel.ActiveWorkbook

# --- xuat-linhkien: drop the obsolete export row first -----------------
$wsXuat = $wb.Worksheets.Item("xuat-linhkien")
$wsXuat.Rows.Item(2).Delete()

# --- nhap-linhkien: refresh the item info for row 2 ---------------------
$wsNhap = $wb.Worksheets.Item("nhap-linhkien")
$wsNhap.Range("A2").Value = "DCU_2G(Sim800C)_4G(Sim7600CE)_SF80P20_Cover_1603757518"
$wsNhap.Range("B2").Value = "DCU-NEW-V4.2-231020"
$wsNhap.Range("D2").Value = "DCU RF 2G"
$wsNhap.Range("H2").Value = 49
